$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Каналы")

# Clear out the old tail of the column-K notes list (rows 5-10 and 13)
# so we can rewrite it with the shifted/updated entries.
$ws.Range("K5:K13").ClearContents()

# Row 3 text is updated (same cell, same bold style, new wording)
$ws.Range("K3").Value = "сделать fixture.channel.count() из DB"

# New row inserted at K4, bold like K2/K3
$ws.Range("K4").Value = "конвертер xls to json"
$ws.Range("K4").Font.Bold = $true

# Existing notes shift down one row, values unchanged
$ws.Range("K6").Value = "Магия с фикстурой pytest: 18 с 08:28"
$ws.Range("K7").Value = "Магия pytest_generate_tests: 57 с 3:20"
$ws.Range("K8").Value = "лекции про строки: 39-44"
$ws.Range("K9").Value = "тесты со случайными данными: 45"
$ws.Range("K10").Value = "DDT: 45, "

# New entries appended after the shifted block
$ws.Range("K11").Value = "параметры запуска в командной строке: 46"
$ws.Range("K12").Value = "json: 50, 56, 58"

# The long-standing note moves far down to K20, with a new note added at K21
$ws.Range("K20").Value = "чтобы тесты запускались из консоли, необоходимо в PATH положить драйвера chromedriver.exe, IEDriverServer.exe"
$ws.Range("K21").Value = "в env проекта должны быть: putest, webdriver, jsonpickle"

$ws.Range("K22").Select()
